# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers table updates
$ws.Range("C3").Value = 4277
$ws.Range("D3").Value = 69.90000000000001

$ws.Range("C4").Value = 6219
$ws.Range("D4").Value = 94.09999999999999

$ws.Range("C5").Value = 10496

# Good Drivers table updates
$ws.Range("B13").Value = 449371
$ws.Range("B14").Value = 77999
